# Dataset_HargaEmas_2025.xlsx - penambahan modul help dan update dataset
# Appends 8 new daily rows (2025-11-01 .. 2025-11-08) to the "Data_Harian_Lengkap"
# sheet, right after the existing last row (305), reusing the existing
# date/percentage number formats from the row above instead of creating new
# style entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: Date (serial), Gold_Price, USD_Sell_Rate, USD_Buy_Rate, BI_Rate
$newData = @(
    @(45962, 2290000, 16723, 16556, 0.0475),
    @(45963, 2290000, 16723, 16556, 0.0475),
    @(45964, 2278000, 16708, 16541, 0.0475),
    @(45965, 2286000, 16747, 16580, 0.0475),
    @(45966, 2260000, 16807, 16640, 0.0475),
    @(45967, 2287000, 16812, 16645, 0.0475),
    @(45968, 2296000, 16790, 16623, 0.0475),
    @(45969, 2299000, 16790, 16623, 0.0475)
)

$startRow = 306

# Copy formatting (number formats / styles) down from the last populated
# row (305) onto the new rows before filling in the values, so the new
# cells reuse the existing "yyyy-mm-dd" and "0.00%"-style formats instead
# of Excel minting brand-new ones.
$lastRow = $startRow + $newData.Count - 1
$ws.Range("A305:E305").Copy()
$ws.Range("A" + $startRow + ":E" + $lastRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]
    $ws.Range("A" + $r).Value = $row[0]
    $ws.Range("B" + $r).Value = $row[1]
    $ws.Range("C" + $r).Value = $row[2]
    $ws.Range("D" + $r).Value = $row[3]
    $ws.Range("E" + $r).Value = $row[4]
}

# Match the updated view/selection state from the edit.
$ws.Activate()
$ws.Range("D307").Select()

Write-Output "Added rows $startRow..$lastRow to $($ws.Name)"
